$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("F41")
# Pre-set look matching hyperlink style but with family=1 Cambria font (fontId 12 u,Cambria,blue,family2 exists already... let's try differently)
$c.Font.Name = "Cambria"
$c.Font.Underline = -4142
$c.Font.Color = 5597121
$c.Value = "tmp"
Write-Host "pre-add done"
$ws.Hyperlinks.Add($c, "https://youtu.be/nJCTSjW0hQM")
Write-Host "post-add font:" $c.Font.Name "underline:" $c.Font.Underline "color:" $c.Font.Color
